$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1300.5
$ws.Range("I43").Value = 1300.5
$ws.Range("K43").Value = 1300.5
$ws.Range("M43").Value = -1231.5
$ws.Range("H51").Value = 6355.15
$ws.Range("I51").Value = 3954.1667
$ws.Range("J51").Value = 7384.143
$ws.Range("K51").Value = 3954.1667
$ws.Range("L51").Value = 7384.143
$ws.Range("M51").Value = -3470.1667
$ws.Range("N51").Value = -8352.143
$ws.Range("H74").Value = 10893.19
$ws.Range("I74").Value = 6595.8887
$ws.Range("K74").Value = 6595.8887
$ws.Range("M74").Value = -5659.8887
$ws.Range("H77").Value = 10893.19
$ws.Range("I77").Value = 6595.8887
$ws.Range("K77").Value = 32979.4435
$ws.Range("M77").Value = -28299.4435
$ws.Range("H138").Value = 7364.8076
$ws.Range("J138").Value = 7484.7236
$ws.Range("L138").Value = 22454.1708
$ws.Range("N138").Value = -32734.1708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3276.2104
$ws.Range("I2").Value = 3095.6667
$ws.Range("K2").Value = 3095.6667
$ws.Range("M2").Value = -2982.6667
$ws.Range("H25").Value = 5762.8
$ws.Range("I25").Value = 4938.3335
$ws.Range("K25").Value = 4938.3335
$ws.Range("M25").Value = -4536.3335
$ws.Range("H35").Value = 7777
$ws.Range("I35").Value = 5554
$ws.Range("K35").Value = 5554
$ws.Range("M35").Value = -5148
$ws.Range("H61").Value = 4960.4
$ws.Range("I61").Value = 3184.9
$ws.Range("K61").Value = 3184.9
$ws.Range("M61").Value = -2972.9
$ws.Range("H102").Value = 305885.34
$ws.Range("I102").Value = 347559.3
$ws.Range("K102").Value = 347559.3
$ws.Range("M102").Value = -345937.3
$ws.Range("H116").Value = 3276.2104
$ws.Range("I116").Value = 3095.6667
$ws.Range("K116").Value = 3095.6667
$ws.Range("M116").Value = -801.6667000000002
$ws.Range("H122").Value = 5141.7896
$ws.Range("J122").Value = 2066.3333
$ws.Range("L122").Value = 6198.999899999999
$ws.Range("N122").Value = -11098.9999
$ws.Range("H136").Value = 4960.4
$ws.Range("I136").Value = 3184.9
$ws.Range("K136").Value = 9554.700000000001
$ws.Range("M136").Value = -7004.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3276.2104
$ws.Range("I3").Value = 3095.6667
$ws.Range("K3").Value = 3095.6667
$ws.Range("M3").Value = -2981.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4396.648
$ws.Range("I31").Value = 3385.6
$ws.Range("J31").Value = 4793.137
$ws.Range("K31").Value = 3385.6
$ws.Range("L31").Value = 4793.137
$ws.Range("M31").Value = -3090.6
$ws.Range("N31").Value = -5383.137
$ws.Range("H34").Value = 4396.648
$ws.Range("I34").Value = 3385.6
$ws.Range("J34").Value = 4793.137
$ws.Range("K34").Value = 3385.6
$ws.Range("L34").Value = 4793.137
$ws.Range("M34").Value = -3183.6
$ws.Range("N34").Value = -5197.137
$ws.Range("H86").Value = 5108.773
$ws.Range("I86").Value = 3674.1333
$ws.Range("J86").Value = 8183
$ws.Range("K86").Value = 3674.1333
$ws.Range("L86").Value = 8183
$ws.Range("M86").Value = -2551.1333
$ws.Range("N86").Value = -10429
$ws.Range("H89").Value = 5108.773
$ws.Range("I89").Value = 3674.1333
$ws.Range("J89").Value = 8183
$ws.Range("K89").Value = 18370.6665
$ws.Range("L89").Value = 40915
$ws.Range("M89").Value = -12754.6665
$ws.Range("N89").Value = -52147
$ws.Range("H103").Value = 51479.266
$ws.Range("I103").Value = 7503
$ws.Range("J103").Value = 67470.63
$ws.Range("K103").Value = 7503
$ws.Range("L103").Value = 67470.63
$ws.Range("M103").Value = -6331
$ws.Range("N103").Value = -69814.63
$ws.Range("H132").Value = 7718.294
$ws.Range("I132").Value = 2593.5557
$ws.Range("K132").Value = 7780.6671
$ws.Range("M132").Value = -5250.6671
$ws.Range("H134").Value = 2891.3225
$ws.Range("I134").Value = 2633.5
$ws.Range("K134").Value = 7900.5
$ws.Range("M134").Value = -5365.5
$ws.Range("H141").Value = 336795.72
$ws.Range("J141").Value = 380952.66
$ws.Range("L141").Value = 380952.66
$ws.Range("N141").Value = -391312.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 121415.75
$ws.Range("I128").Value = 121415.75
$ws.Range("K128").Value = 364247.25
$ws.Range("M128").Value = -359267.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8993.929
$ws.Range("I113").Value = 10665.091
$ws.Range("K113").Value = 10665.091
$ws.Range("M113").Value = -8495.091
$ws.Range("H132").Value = 6604.353
$ws.Range("I132").Value = 5559.769
$ws.Range("J132").Value = 9999.25
$ws.Range("K132").Value = 16679.307
$ws.Range("L132").Value = 29997.75
$ws.Range("M132").Value = -14149.307
$ws.Range("N132").Value = -35057.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 143030560
$ws.Range("I7").Value = 143030560
$ws.Range("K7").Value = 143030560
$ws.Range("M7").Value = -143030448
$ws.Range("H22").Value = 837143.75
$ws.Range("I22").Value = 1574.25
$ws.Range("J22").Value = 1254928.5
$ws.Range("K22").Value = 1574.25
$ws.Range("L22").Value = 1254928.5
$ws.Range("M22").Value = -1279.25
$ws.Range("N22").Value = -1255518.5
$ws.Range("H27").Value = 837143.75
$ws.Range("I27").Value = 1574.25
$ws.Range("J27").Value = 1254928.5
$ws.Range("K27").Value = 1574.25
$ws.Range("L27").Value = 1254928.5
$ws.Range("M27").Value = -1467.25
$ws.Range("N27").Value = -1255142.5
$ws.Range("H40").Value = 63653.816
$ws.Range("I40").Value = 77385.03
$ws.Range("K40").Value = 77385.03
$ws.Range("M40").Value = -77249.03
$ws.Range("H126").Value = 143030560
$ws.Range("I126").Value = 143030560
$ws.Range("K126").Value = 429091680
$ws.Range("M126").Value = -429089210
$ws.Range("H132").Value = 5138.5864
$ws.Range("I132").Value = 4015.652
$ws.Range("K132").Value = 12046.956
$ws.Range("M132").Value = -9516.956
$ws.Range("H136").Value = 5042.385
$ws.Range("I136").Value = 3858.6843
$ws.Range("J136").Value = 8255.286
$ws.Range("K136").Value = 11576.0529
$ws.Range("L136").Value = 24765.858
$ws.Range("M136").Value = -9026.052899999999
$ws.Range("N136").Value = -29865.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7920.3335
$ws.Range("I132").Value = 7994.5386
$ws.Range("K132").Value = 23983.6158
$ws.Range("M132").Value = -21453.6158
$ws.Range("H136").Value = 20004774
$ws.Range("I136").Value = 24395434
$ws.Range("J136").Value = 2884.2222
$ws.Range("K136").Value = 73186302
$ws.Range("L136").Value = 8652.6666
$ws.Range("M136").Value = -73183752
$ws.Range("N136").Value = -13752.6666
